$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.146172523498535
$ws.Range("B1").Value = 2.79360294342041
$ws.Range("C1").Value = 6.898158550262451
$ws.Range("D1").Value = 1.989449501037598
$ws.Range("E1").Value = 1.047636985778809
